# Automatic update of files.
#
# 1) Bump the "Förändrad" (changed) date in column C, rows 2-18, from
#    2023-10-08 (serial 45207) to 2023-10-09 (serial 45208).
# 2) Repoint the hyperlink formulas in row 2 (columns S-Y) from the
#    "Logging_VASTERVIK" folder to "Logging_0883".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column C: rows 2 through 18 ---
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}

# --- 2) Row 2, columns S through Y: swap the folder name in the URL ---
$cols = @("S", "T", "U", "V", "W", "X", "Y")
foreach ($col in $cols) {
    $rng = $ws.Range($col + "2")
    $formula = $rng.Formula
    if ($formula -and $formula.Contains("Logging_VASTERVIK")) {
        $rng.Formula = $formula.Replace("Logging_VASTERVIK", "Logging_0883")
    }
}
